$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C2:C6 to be formulas referencing D column divided by 3
$ws.Range("C2").Formula = "=D2/3"
$ws.Range("C3:C6").Formula = "=D3/3"

# Update the active selection to E10
$ws.Range("E10").Select()
